$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 7 (row 8) is now "In process": fill in the Status cell F8, reusing the
# same look already applied to the other "In process" rows (F3/F5/F7/F10 -
# yellow fill, no text wrap) by copying their format across.
$ws.Range("F3").Copy()
$ws.Range("F8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F8").Value = "In process"

# Leave the newly-edited cell selected.
$ws.Range("F8").Select()
